$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# File 26275d01-d2ec-4257-8b7a-e8bc11d8ad91.md was just handed off again -
# refresh its handoff/generate timestamps across all report sheets.
$wsOverview.Range("G6").Value = "2016-08-21 04:48:34"
$wsZhCn.Range("H6").Value = "2016-08-21 04:48:30"
$wsDeDe.Range("H6").Value = "2016-08-21 04:48:34"
